$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 369, pushing the existing row 369 (and all
# rows below it, down through the former row 402) down by one row. This
# turns the former rows 369..402 into the new rows 370..403, matching the
# diff (which shows every one of those rows' data shifting down by one
# position, with the last original row 402 ending up duplicated into the
# new row 403).
$ws.Rows.Item(369).Insert()

# Populate the newly inserted row 369 with this week's new record.
$ws.Range("A369").Value = 7
$ws.Range("B369").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C369").Value = "Ñuble"
$ws.Range("D369").Value = 45013
$ws.Range("E369").Value = 16
$ws.Range("F369").Value = 100114013
$ws.Range("G369").Value = "Zanahoria"
$ws.Range("H369").Value = "Sin especificar"
$ws.Range("I369").Value = "Primera"
$ws.Range("J369").Value = 150
$ws.Range("K369").Value = 7000
$ws.Range("L369").Value = 7000
$ws.Range("M369").Value = 7000
$ws.Range("N369").Value = "$/saco 20 kilos"
$ws.Range("O369").Value = "Provincia de Diguillín"
$ws.Range("P369").Value = 350
$ws.Range("Q369").Value = 20
$ws.Range("R369").Value = "Hortaliza"
